$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 109 (shifts old 109-112 down to 110-113)
$ws.Rows.Item(109).Insert()

# 2. Populate the new row 109: "* Registration" / 2 hours, matching the
#    formatting of the rows above it (B/C use the bordered "entry" style,
#    E gets the blank "tail" style used elsewhere in the sheet, e.g. E43/E44)
$ws.Range("B108").Copy()
$ws.Range("B109").PasteSpecial(-4122)
$ws.Range("C108").Copy()
$ws.Range("C109").PasteSpecial(-4122)
$ws.Range("E43").Copy()
$ws.Range("E109").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B109").Value = "* Registration"
$ws.Range("C109").Value = 2

# 3. C105 hours 7 -> 9
$ws.Range("C105").Value = 9

# 4. Update selection to match the authored state
$null = $ws.Range("C111").Select()
